$wb = $excel.ActiveWorkbook

# --- "Performance Metrics" sheet: tidy up confidence-interval spacing ---
$wsPM = $wb.Worksheets.Item("Performance Metrics")
$wsPM.Range("J4").Value = "1.24 [1.2,1.29]"
$wsPM.Range("K2").Value = "1.53 [1.5,1.56]"
$wsPM.Range("K3").Value = "1.64 [1.6,1.68]"
$wsPM.Range("N2").Value = "0.522 [0.519,0.527]"
$wsPM.Range("N4").Value = "0.581 [0.571,0.592]"

# --- "Evaluation Sample Sets" sheet: switch dict-like text to REST API v1.4 format ---
$wsESS = $wb.Worksheets.Item("Evaluation Sample Sets")
$wsESS.Range("G4").Value = "mean:54.3;range:[50.1,58.4];unit:years"
$wsESS.Range("M4").Value = "mean:6.2;sd:1.7;unit:years"
